$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new progress note for Pradnya Kore (row 9)
$ws.Range("B9").Value = "studied basics of dart and setup of flutter env"

# Update the active selection to match the edited cell
$ws.Range("B9").Select()
